# Update "想去人数" (want-to-go count) values in column F across sheets,
# per the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 1239
$ws1.Range("F10").Value = 348
$ws1.Range("F12").Value = 2484
$ws1.Range("F20").Value = 405222
$ws1.Range("F21").Value = 1204
$ws1.Range("F23").Value = 587
$ws1.Range("F29").Value = 172
$ws1.Range("F30").Value = 1239
$ws1.Range("F31").Value = 409
$ws1.Range("F32").Value = 129
$ws1.Range("F41").Value = 123
$ws1.Range("F43").Value = 2662

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F7").Value  = 169
$ws3.Range("F10").Value = 894
$ws3.Range("F13").Value = 1298
$ws3.Range("F15").Value = 1067

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 169
$ws4.Range("F8").Value  = 894
$ws4.Range("F12").Value = 1298
$ws4.Range("F16").Value = 1239
$ws4.Range("F17").Value = 348
$ws4.Range("F18").Value = 1067
$ws4.Range("F19").Value = 2484
$ws4.Range("F28").Value = 1204
$ws4.Range("F34").Value = 172
$ws4.Range("F37").Value = 1239
$ws4.Range("F38").Value = 409
$ws4.Range("F46").Value = 123
$ws4.Range("F48").Value = 2662
